$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple: row number, Vencimento (A), Taxa (B), Data de Salvamento (C)
$rows = @(
    ,@(2, 45792, 0, "2025-04-04 13:05:55")
    ,@(3, 45792, 1, "2025-04-04 13:06:11")
    ,@(4, 45792, 3, "2025-04-04 13:06:27")
    ,@(5, 46249, 0, "2025-04-04 13:05:55")
    ,@(6, 46249, 1, "2025-04-04 13:06:11")
    ,@(7, 46249, 3, "2025-04-04 13:06:27")
    ,@(8, 46522, 3, "2025-04-04 13:06:27")
    ,@(9, 46522, 0, "2025-04-04 13:05:55")
    ,@(10, 46522, 1, "2025-04-04 13:06:11")
    ,@(11, 46980, 1, "2025-04-04 13:06:11")
    ,@(12, 46980, 3, "2025-04-04 13:06:27")
    ,@(13, 46980, 0, "2025-04-04 13:05:55")
    ,@(14, 47253, 0, "2025-04-04 13:05:55")
    ,@(15, 47253, 1, "2025-04-04 13:06:11")
    ,@(16, 47253, 3, "2025-04-04 13:06:27")
    ,@(17, 47710, 0, "2025-04-04 13:05:55")
    ,@(18, 47710, 1, "2025-04-04 13:06:11")
    ,@(19, 47710, 3, "2025-04-04 13:06:27")
    ,@(20, 48441, 0, "2025-04-04 13:05:55")
    ,@(21, 48441, 1, "2025-04-04 13:06:11")
    ,@(22, 48441, 3, "2025-04-04 13:06:27")
    ,@(23, 48714, 3, "2025-04-04 13:06:27")
    ,@(24, 48714, 0, "2025-04-04 13:05:55")
    ,@(25, 48714, 1, "2025-04-04 13:06:11")
    ,@(26, 49444, 1, "2025-04-04 13:06:11")
    ,@(27, 49444, 0, "2025-04-04 13:05:55")
    ,@(28, 49444, 3, "2025-04-04 13:06:27")
    ,@(29, 51363, 3, "2025-04-04 13:06:27")
    ,@(30, 51363, 0, "2025-04-04 13:05:55")
    ,@(31, 51363, 1, "2025-04-04 13:06:11")
    ,@(32, 53097, 3, "2025-04-04 13:06:27")
    ,@(33, 53097, 1, "2025-04-04 13:06:11")
    ,@(34, 53097, 0, "2025-04-04 13:05:55")
    ,@(35, 55015, 0, "2025-04-04 13:05:55")
    ,@(36, 55015, 1, "2025-04-04 13:06:11")
    ,@(37, 55015, 3, "2025-04-04 13:06:27")
    ,@(38, 56749, 0, "2025-04-04 13:05:55")
    ,@(39, 56749, 1, "2025-04-04 13:06:11")
    ,@(40, 56749, 3, "2025-04-04 13:06:27")
    ,@(41, 58668, 0, "2025-04-04 13:05:55")
    ,@(42, 58668, 1, "2025-04-04 13:06:11")
    ,@(43, 58668, 3, "2025-04-04 13:06:27")
)

foreach ($item in $rows) {
    $r = $item[0]
    $ws.Cells.Item($r, 1).Value = $item[1]
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 2).Value = $item[2]
    $ws.Cells.Item($r, 3).Value = $item[3]
}

$ws.Range("A1").Select() | Out-Null
